$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.Calculation = "xlCalculationManual"

Write-Output "done"
